$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.111.82'
$ws.Range('E2').Value = '  +0.02%  '
$ws.Range('D3').Value = '1.786.52'
$ws.Range('E3').Value = '  -0.26%  '
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '226.62'
$ws.Range('E5').Value = '  -0.81%  '
$ws.Range('E6').Value = '  -0.82%  '
$ws.Range('E7').Value = '  +0.26%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '31.95'
$ws.Range('E8').Value = '  -2.47%  '
$ws.Range('E9').Value = '  +0.82%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0691'
$ws.Range('E10').Value = '  -3.32%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0944'
$ws.Range('E11').Value = '  +0.75%  '
$ws.Range('D12').Value = '2.043.86'
$ws.Range('E12').Value = '  -0.26%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.23'
$ws.Range('E13').Value = '  +0.78%  '
$ws.Range('D14').Value = '1.791.17'
$ws.Range('E14').Value = '  -0.22%  '
$ws.Range('D15').Value = '34.054.61'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.620'
$ws.Range('E16').Value = '  -0.95%  '
$ws.Range('E17').Value = '  +0.60%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '68.00'
$ws.Range('E18').Value = '  -0.74%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '245.75'
$ws.Range('E19').Value = '  +0.11%  '
$ws.Range('E20').Value = '  -1.22%  '
$ws.Range('E21').Value = '  +0.20%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.82'
$ws.Range('E22').Value = '  +0.44%  '
$ws.Range('E23').Value = '  -0.58%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.04'
$ws.Range('E24').Value = '  -2.25%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '161.61'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.15'
$ws.Range('E26').Value = '  +0.33%  '
$ws.Range('E27').Value = '  -0.27%  '
$ws.Range('E28').Value = '  +0.46%  '
$ws.Range('E29').Value = '  +0.36%  '
$ws.Range('E30').Value = '  -0.90%  '
$ws.Range('E31').Value = '  +0.45%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.66'
$ws.Range('E32').Value = '  -0.40%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.61'
$ws.Range('E33').Value = '  +2.50%  '
$ws.Range('E34').Value = '  -1.02%  '
$ws.Range('D35').Value = '1.449.51'
$ws.Range('E35').Value = '  +3.43%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.645'
$ws.Range('E36').Value = '  -2.16%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0193'
$ws.Range('E37').Value = '  +2.17%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.39'
$ws.Range('E38').Value = '  +7.85%  '
$ws.Range('E39').Value = '  -1.01%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '80.15'
$ws.Range('E40').Value = '  +1.60%  '
$ws.Range('E41').Value = '  +0.37%  '
$ws.Range('E42').Value = '  -0.02%  '
$ws.Range('E43').Value = '  -0.51%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '13.55'
$ws.Range('E44').Value = '  +1.80%  '
$ws.Range('E45').Value = '  +1.96%  '
$ws.Range('E46').Value = '  +3.36%  '
$ws.Range('E47').Value = '  -0.34%  '
$ws.Range('E48').Value = '  -0.77%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '107.58'
$ws.Range('E49').Value = '  -1.57%  '
$ws.Range('D50').Value = '1.945.39'
$ws.Range('E50').Value = '  -0.11%  '
$ws.Range('E51').Value = '  +0.22%  '
